$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.423.64'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.642.71'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.99'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.531'
$ws.Range('E6').Value = '  +4.07%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.09'
$ws.Range('E8').Value = '  -3.41%  '
$ws.Range('E9').Value = '  -2.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0611'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Value = '1.875.01'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').Value = '1.642.13'
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.570'
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.44'
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').Value = '27.393.24'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.05'
$ws.Range('E18').Value = '  -5.29%  '
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.58'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  -3.66%  '
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.37'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.01'
$ws.Range('E26').Value = '  -2.91%  '
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.60'
$ws.Range('E29').Value = '  -4.94%  '
$ws.Range('E30').Value = '  -3.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0484'
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.29'
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('D34').Value = '1.416.20'
$ws.Range('E34').Value = '  -3.12%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.884'
$ws.Range('E38').Value = '  -4.24%  '
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.825'
$ws.Range('E42').Value = '  +4.76%  '
$ws.Range('B43').Value = 'mCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.45'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.50'
$ws.Range('E44').Value = '  +1.64%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.24'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.60'
$ws.Range('E46').Value = '  -7.21%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.784.26'
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.67'
$ws.Range('E48').Value = '  -3.96%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.03'
$ws.Range('E49').Value = '  -1.20%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0991'
$ws.Range('E51').Value = '  -3.27%  '
